# Update retention metrics data following a data refresh.
# - Row 31 (cohort 2023, period_index 3): num_customers 39 -> 40, retention_rate recalculated
# - Row 36 (cohort 2024, period_index 1): num_customers 107 -> 109, retention_rate recalculated
# - Row 37 (cohort 2025, period_index 0): num_customers 654 -> 664, cohort_size 654 -> 664

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: cohort 2023, period 3
$ws.Range("C31").Value = 40
$ws.Range("E31").Value = 0.01730103806228374

# Row 36: cohort 2024, period 1
$ws.Range("C36").Value = 109
$ws.Range("E36").Value = 0.05647668393782383

# Row 37: cohort 2025, period 0
$ws.Range("C37").Value = 664
$ws.Range("D37").Value = 664
